# Add a new "November" worksheet after the existing "August" sheet,
# populate it like the "Template" sheet (same labels) but with the
# MONTH/NAME fields left as "xxx" placeholders, and give it the same
# page margins used on the other sheets.

$wb = $excel.ActiveWorkbook

# Insert the new sheet right after the last existing sheet ("August")
# so the final order is Template, August, November.
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add([System.Type]::Missing, $lastSheet)
$ws.Name = "November"

# Match the page margins used elsewhere in the workbook
# (0.75in left/right, 1in top/bottom, 0.5in header/footer).
$ws.PageSetup.LeftMargin = 54
$ws.PageSetup.RightMargin = 54
$ws.PageSetup.TopMargin = 72
$ws.PageSetup.BottomMargin = 72
$ws.PageSetup.HeaderMargin = 36
$ws.PageSetup.FooterMargin = 36

# Header row
$ws.Range("A1").Value = "TIMESHEET"
$ws.Range("B1").Value = "CNR"
$ws.Range("C1").Value = ".Architects"

# Month / name placeholders
$ws.Range("A3").Value = "MONTH:"
$ws.Range("B3").Value = "xxx"

$ws.Range("A4").Value = "NAME:"
$ws.Range("B4").Value = "xxx"

# Table header row
$ws.Range("A6").Value = "PROJECT"
$ws.Range("B6").Value = "DESCRIPTION OF WORK"
$ws.Range("C6").Value = "OFFICE USE"
